# Update the cryptos list with latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $value) {
    # Force the cell to stay a text value even when the new string looks
    # like a plain number (e.g. "613.48"), matching the scraper's original
    # inline-string output instead of Excel's automatic numeric coercion.
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.236.00"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.725.16"
$ws.Range("E3").Value = "  -0.44%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "613.48"
$ws.Range("E5").Value = "  +5.14%  "

# Row 6 - Solana
Set-TextValue "D6" "191.51"
$ws.Range("E6").Value = "  +7.96%  "

# Row 7 - XRP
Set-TextValue "D7" "0.637"
$ws.Range("E7").Value = "  +0.29%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.20%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.725"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10 - was Dogecoin, now Avalanche
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D10" "60.74"
$ws.Range("E10").Value = "  +12.28%  "

# Row 11 - was Avalanche, now Dogecoin
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.162"
$ws.Range("E11").Value = "  -3.24%  "

# Row 12 - ShibaInu
Set-TextValue "D12" "0.0000291"
$ws.Range("E12").Value = "  -3.55%  "

# Row 13 - Polkadot
Set-TextValue "D13" "10.65"
$ws.Range("E13").Value = "  -1.43%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.323.11"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15 - WrappedEther
Set-TextValue "D15" "3.722.03"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16 - Chainlink
Set-TextValue "D16" "19.43"
$ws.Range("E16").Value = "  -0.56%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -0.42%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.22%  "

# Row 19 - Uniswap
Set-TextValue "D19" "12.94"
$ws.Range("E19").Value = "  -1.74%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "69.094.26"
$ws.Range("E20").Value = "  +1.14%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "412.40"
$ws.Range("E21").Value = "  -0.21%  "

# Row 22 - PancakeSwap
$ws.Range("E22").Value = "  +0.32%  "

# Row 23 - Litecoin
Set-TextValue "D23" "89.66"
$ws.Range("E23").Value = "  +0.46%  "

# Row 24 - ImmutableX
Set-TextValue "D24" "3.06"
$ws.Range("E24").Value = "  -1.50%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "12.88"
$ws.Range("E25").Value = "  -0.74%  "

# Row 26 - RenderToken
Set-TextValue "D26" "10.89"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27 - Toncoin
Set-TextValue "D27" "3.81"
$ws.Range("E27").Value = "  -1.76%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  +1.03%  "

# Row 29 - Filecoin
Set-TextValue "D29" "9.71"
$ws.Range("E29").Value = "  +0.72%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "33.15"
$ws.Range("E30").Value = "  -0.35%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.73"
$ws.Range("E31").Value = "  -3.84%  "

# Row 32 - Cosmos
Set-TextValue "D32" "12.82"
$ws.Range("E32").Value = "  +0.05%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.124"
$ws.Range("E33").Value = "  +4.02%  "

# Row 34 - InjectiveProtocol
Set-TextValue "D34" "45.97"
$ws.Range("E34").Value = "  +4.16%  "

# Row 35 - Bittensor
Set-TextValue "D35" "631.71"
$ws.Range("E35").Value = "  +2.50%  "

# Row 36 - OKB
$ws.Range("E36").Value = "  -0.70%  "

# Row 37 - was PEPE, now TheGraph
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D37" "0.415"
$ws.Range("E37").Value = "  +2.96%  "

# Row 38 - was TheGraph, now PEPE
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0830"
$ws.Range("E38").Value = "  -10.60%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  -0.13%  "

# Row 40 - FirstDigitalUSD
Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  +0.13%  "

# Row 41 - Kaspa
Set-TextValue "D41" "0.142"
$ws.Range("E41").Value = "  +2.92%  "

# Row 42 - ThetaToken
$ws.Range("E42").Value = "  -1.38%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -0.19%  "

# Row 44 - Fetch.AI
$ws.Range("E44").Value = "  -0.08%  "

# Row 45 - Stellar
$ws.Range("E45").Value = "  +2.72%  "

# Row 46 - Maker
Set-TextValue "D46" "2.885.05"
$ws.Range("E46").Value = "  +5.46%  "

# Row 47 - THORChain
$ws.Range("E47").Value = "  -2.81%  "

# Row 48 - WEMIXToken
$ws.Range("E48").Value = "  +1.19%  "

# Row 49 - ApeXProtocol
Set-TextValue "D49" "3.11"
$ws.Range("E49").Value = "  -1.66%  "

# Row 50 - Monero
Set-TextValue "D50" "142.31"
$ws.Range("E50").Value = "  -1.06%  "

# Row 51 - was dogwifhat, now Stacks
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D51" "2.79"
$ws.Range("E51").Value = "  +0.02%  "
